# Documento Requisitos (Excel) atualizado
# Applies the requirements-sheet update described by the commit diff:
#   - R05 description is amended and renumbered to R05.1 (row 14)
#   - C13 gains center/center alignment formatting
#   - Four new requirement rows (R10..R13) are filled in (rows 19-22)
#   - Selection / scroll position move to reflect the newly entered rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documento Geral de Recolha")

# --- Row 14: R05 -> R05.1, updated description + timestamp -----------------
$ws.Range("A14").Value = "R05.1"
$ws.Range("B14").Value = 45560.4375
$ws.Range("C14").Value = "Cada funcionário tem um nome, um número de id e um cargo/função (juíz, camera, auxiliar)"

# --- Row 13: purely a formatting touch-up (center the description cell) ----
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("C13").VerticalAlignment = -4108

# --- Row 19: R10 ------------------------------------------------------------
$ws.Range("A19").Value = "R10"
$ws.Range("B19").Value = "25/09/202410:16"
$ws.Range("C19").Value = "Deve ser possível gerir os funcionários da competição"
$ws.Range("F19").Value = "JF "
$ws.Range("G19").Value = "M"

# --- Row 20: R11 ------------------------------------------------------------
$ws.Range("A20").Value = "R11"
$ws.Range("B20").Value = 45560.428472222222
$ws.Range("C20").Value = "Deve ser possível adicionar, editar e remover atletas e treinadores antes do início da competição"
$ws.Range("F20").Value = "JF"
$ws.Range("G20").Value = "M"

# --- Row 21: R12 ------------------------------------------------------------
$ws.Range("A21").Value = "R12"
$ws.Range("B21").Value = 45560.430555555555
$ws.Range("C21").Value = "Deve ser possível adicionar e editar os resultados das partidas"
$ws.Range("F21").Value = "JF"
$ws.Range("G21").Value = "M"

# --- Row 22: R13 ------------------------------------------------------------
$ws.Range("A22").Value = "R13"
$ws.Range("B22").Value = 45560.438888888886
$ws.Range("C22").Value = "Cada partida acontecerá a uma determinada hora, entre duas ou mais equipes"
$ws.Range("F22").Value = "JF"
$ws.Range("G22").Value = "D"

# --- Reflect the new scroll position / active selection --------------------
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("A23").Select()
